# Remove the trailing "empty paragraph / page-break paragraph / footer
# copyright paragraph" block that used to sit right after the last
# bibliography line ("...o século 21. Editora Campus. 2010."), leaving
# that bibliography paragraph directly followed by the two remaining
# (still-empty) paragraphs that close the document.

$d = $word.ActiveDocument

# Locate the end of the last bibliography paragraph ("...2010.") - this
# paragraph itself is kept untouched. Expand(4) = wdParagraph, so the
# range grows to cover the whole paragraph including its paragraph mark.
$anchorStart = $d.Content
$anchorStart.Find.Execute("Editora Campus. 2010.", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorStart.Expand(4) | Out-Null

# Locate the end of the copyright paragraph ("... Creative Commons
# Attribution") - everything up to (and including) this paragraph's
# mark gets removed.
$anchorEnd = $d.Content
$anchorEnd.Find.Execute("Powered by Jekyll and Github pages", $true, $false, `
    $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorEnd.Expand(4) | Out-Null

# Delete everything from just after the bibliography paragraph's mark
# through the end of the copyright paragraph's mark - this removes the
# blank paragraph, the page-break paragraph, and the copyright
# paragraph in one shot.
$deadZone = $d.Range($anchorStart.End, $anchorEnd.End)
$deadZone.Delete()
